$d = $word.ActiveDocument

# NOTE: this COM-interop shim re-seats any stored Tables.Item()/Rows.Item()/
# Cell() reference to whatever Tables.Item() was most recently resolved
# elsewhere in the script (a caching quirk of this particular host), so
# every mutation below walks the full $d.Tables.Item(n).Cell(row, col) path
# inline right at the point of use instead of caching table/row/cell
# objects in variables across statements.

# --- Header info table (Tables.Item(1)): exam date + CIE number ---
$d.Tables.Item(1).Cell(1, 4).Range.Text = "2015/11/19"   # Date : 2015/11/11 -> 2015/11/19
$d.Tables.Item(1).Cell(3, 4).Range.Text = "      3"       # CIE  :       2   ->       3

# --- Question table A (Tables.Item(2)): Sl No 1 and 2 ---
# Sl No 1: question text + marks
$d.Tables.Item(2).Cell(2, 2).Range.Text = "This is Question number 321"
$d.Tables.Item(2).Cell(2, 3).Range.Text = "1"
# Sl No 2: was fully blank, now filled in (question, marks, CO, LO)
$d.Tables.Item(2).Cell(3, 2).Range.Text = "This is Question number 977"
$d.Tables.Item(2).Cell(3, 3).Range.Text = "2"
$d.Tables.Item(2).Cell(3, 4).Range.Text = "CO-1"
$d.Tables.Item(2).Cell(3, 5).Range.Text = "LO-1"

# --- Question table B (Tables.Item(3)): Sl No 1a..5c ---
# 1a
$d.Tables.Item(3).Cell(2, 2).Range.Text = "This is Question number 193"
$d.Tables.Item(3).Cell(2, 3).Range.Text = "3"
# 1b
$d.Tables.Item(3).Cell(3, 2).Range.Text = "This is Question number 209"
$d.Tables.Item(3).Cell(3, 3).Range.Text = "4"
# 1c
$d.Tables.Item(3).Cell(4, 2).Range.Text = "This is Question number 353"
$d.Tables.Item(3).Cell(4, 3).Range.Text = "3"
# 2a
$d.Tables.Item(3).Cell(5, 2).Range.Text = "This is Question number 1025"
$d.Tables.Item(3).Cell(5, 3).Range.Text = "5"
# 2b
$d.Tables.Item(3).Cell(6, 2).Range.Text = "This is Question number 1489"
$d.Tables.Item(3).Cell(6, 3).Range.Text = "4"
# 2c (question text only; marks unchanged in the source diff)
$d.Tables.Item(3).Cell(7, 2).Range.Text = "This is Question number 1761"
# 3a
$d.Tables.Item(3).Cell(8, 2).Range.Text = "This is Question number 2113"
$d.Tables.Item(3).Cell(8, 3).Range.Text = "3"
# 3b
$d.Tables.Item(3).Cell(9, 2).Range.Text = "This is Question number 2449"
$d.Tables.Item(3).Cell(9, 3).Range.Text = "4"
# 3c
$d.Tables.Item(3).Cell(10, 2).Range.Text = "This is Question number 2753"
$d.Tables.Item(3).Cell(10, 3).Range.Text = "3"
# 4a
$d.Tables.Item(3).Cell(11, 2).Range.Text = "This is Question number 2961"
$d.Tables.Item(3).Cell(11, 3).Range.Text = "6"
# 4b
$d.Tables.Item(3).Cell(12, 2).Range.Text = "This is Question number 353"
$d.Tables.Item(3).Cell(12, 3).Range.Text = "3"
# 4c (question text only; marks unchanged in the source diff)
$d.Tables.Item(3).Cell(13, 2).Range.Text = "This is Question number 641"
# 5a
$d.Tables.Item(3).Cell(14, 2).Range.Text = "This is Question number 1137"
$d.Tables.Item(3).Cell(14, 3).Range.Text = "2"
# 5b
$d.Tables.Item(3).Cell(15, 2).Range.Text = "This is Question number 1505"
$d.Tables.Item(3).Cell(15, 3).Range.Text = "5"
# 5c
$d.Tables.Item(3).Cell(16, 2).Range.Text = "This is Question number 1793"
$d.Tables.Item(3).Cell(16, 3).Range.Text = "3"

Write-Host "Applied 36 cell edits"
